$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 79; existing rows 79..158 shift down to 80..159
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly record
$ws.Cells.Item(79, 1).Value = 1
$ws.Cells.Item(79, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(79, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(79, 4).Value = 44981
$ws.Cells.Item(79, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(79, 5).Value = 15
$ws.Cells.Item(79, 6).Value = 100112042
$ws.Cells.Item(79, 7).Value = "Locoto"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 60
$ws.Cells.Item(79, 11).Value = 54000
$ws.Cells.Item(79, 12).Value = 55000
$ws.Cells.Item(79, 13).Value = 54500
$ws.Cells.Item(79, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(79, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(79, 16).Value = 2725
$ws.Cells.Item(79, 17).Value = 20
$ws.Cells.Item(79, 18).Value = "Hortaliza"
